$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.065281189239591
$ws.Range("D2").Value = 1.063798792257822
$ws.Range("E2").Value = 1.06898625276611
$ws.Range("F2").Value = 1.073387369254169
$ws.Range("I2").Value = 1.045818357153772
$ws.Range("J2").Value = 1.070237108864557
$ws.Range("K2").Value = 1.066516456398953
$ws.Range("L2").Value = 1.071689968184904
$ws.Range("M2").Value = 1.076079364656123
$ws.Range("N2").Value = 1.071756968317057

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.067455819005514
$ws.Range("D3").Value = 1.065498133585349
$ws.Range("E3").Value = 1.071114122387056
$ws.Range("F3").Value = 1.075348224262384
$ws.Range("I3").Value = 1.046442479031992
$ws.Range("J3").Value = 1.072061487455489
$ws.Range("K3").Value = 1.068028662588434
$ws.Range("L3").Value = 1.073630669320397
$ws.Range("M3").Value = 1.077854333485588
$ws.Range("N3").Value = 1.073583937734847

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.068855993169475
$ws.Range("D4").Value = 1.066591204971827
$ws.Range("E4").Value = 1.072484425594976
$ws.Range("F4").Value = 1.076610593678994
$ws.Range("I4").Value = 1.046841743645127
$ws.Range("J4").Value = 1.073234900291698
$ws.Range("K4").Value = 1.069000192154893
$ws.Range("L4").Value = 1.074879496306773
$ws.Range("M4").Value = 1.078996017961631
$ws.Range("N4").Value = 1.074759016951875

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.069442998496176
$ws.Range("D5").Value = 1.067049201523925
$ws.Range("E5").Value = 1.073058963823473
$ws.Range("F5").Value = 1.077139786025386
$ws.Range("I5").Value = 1.047008510047101
$ws.Range("J5").Value = 1.073726540052596
$ws.Range("K5").Value = 1.069406982679317
$ws.Range("L5").Value = 1.075402876404691
$ws.Range("M5").Value = 1.079474375893568
$ws.Range("N5").Value = 1.075251354897626

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.069541464896895
$ws.Range("D6").Value = 1.067126012247233
$ws.Range("E6").Value = 1.073155342208812
$ws.Range("F6").Value = 1.077228552278922
$ws.Range("I6").Value = 1.047036447641699
$ws.Range("J6").Value = 1.073808991949898
$ws.Range("K6").Value = 1.069475189202544
$ws.Range("L6").Value = 1.075490659807712
$ws.Range("M6").Value = 1.079554601001741
$ws.Range("N6").Value = 1.07533392388608

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.068863843109595
$ws.Range("D7").Value = 1.066597330719331
$ws.Range("E7").Value = 1.072492108595249
$ws.Range("F7").Value = 1.07661767065356
$ws.Range("I7").Value = 1.046843976229657
$ws.Range("J7").Value = 1.073241476097176
$ws.Range("K7").Value = 1.069005634118493
$ws.Range("L7").Value = 1.074886496084028
$ws.Range("M7").Value = 1.079002416073459
$ws.Range("N7").Value = 1.074765602095752

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.066017578758562
$ws.Range("D8").Value = 1.06437445958767
$ws.Range("E8").Value = 1.069706758491091
$ws.Range("F8").Value = 1.074051401645011
$ws.Range("I8").Value = 1.046030239078116
$ws.Range("J8").Value = 1.070855152242937
$ws.Range("K8").Value = 1.067028973437808
$ws.Range("L8").Value = 1.072347293965112
$ws.Range("M8").Value = 1.076680658589154
$ws.Range("N8").Value = 1.072375889387925

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.060947112338397
$ws.Range("D9").Value = 1.06040624039954
$ws.Range("E9").Value = 1.064746666511745
$ws.Range("F9").Value = 1.069478536381993
$ws.Range("I9").Value = 1.044560632707084
$ws.Range("J9").Value = 1.066594429677077
$ws.Range("K9").Value = 1.063491204743464
$ws.Range("L9").Value = 1.067818234959403
$ws.Range("M9").Value = 1.072535635959192
$ws.Range("N9").Value = 1.068109116107309

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.05752742379537
$ws.Range("D10").Value = 1.057724431209534
$ws.Range("E10").Value = 1.061402714277179
$ws.Range("F10").Value = 1.066393721080516
$ws.Range("I10").Value = 1.04355605921462
$ws.Range("J10").Value = 1.063714404171736
$ws.Range("K10").Value = 1.061094165274638
$ws.Range("L10").Value = 1.064759941320083
$ws.Range("M10").Value = 1.069734129066229
$ws.Range("N10").Value = 1.065225000635407

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.056036777779524
$ws.Range("D11").Value = 1.056554136648532
$ws.Range("E11").Value = 1.059945398943952
$ws.Range("F11").Value = 1.065048892317129
$ws.Range("I11").Value = 1.043114988259216
$ws.Range("J11").Value = 1.062457471788774
$ws.Range("K11").Value = 1.060046678512175
$ws.Range("L11").Value = 1.063425945210272
$ws.Range("M11").Value = 1.06851154507286
$ws.Range("N11").Value = 1.06396628326429

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.055481551117163
$ws.Range("D12").Value = 1.056118039985306
$ws.Range("E12").Value = 1.059402635557801
$ws.Range("F12").Value = 1.064547957079819
$ws.Range("I12").Value = 1.042950224331989
$ws.Range("J12").Value = 1.061989068500613
$ws.Range("K12").Value = 1.059656124884099
$ws.Range("L12").Value = 1.062928934508986
$ws.Range("M12").Value = 1.068055954554584
$ws.Range("N12").Value = 1.063497214789733

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.055600719322273
$ws.Range("D13").Value = 1.056211648025318
$ws.Range("E13").Value = 1.059519126566363
$ws.Range("F13").Value = 1.064655473682305
$ws.Range("I13").Value = 1.042985609130471
$ws.Range("J13").Value = 1.062089612163947
$ws.Range("K13").Value = 1.0597399670743
$ws.Range("L13").Value = 1.063035613795502
$ws.Range("M13").Value = 1.068153747351673
$ws.Range("N13").Value = 1.063597901236601

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.05599091414756
$ws.Range("D14").Value = 1.056518117480348
$ws.Range("E14").Value = 1.059900563828818
$ws.Range("F14").Value = 1.065007513776915
$ws.Range("I14").Value = 1.043101387900788
$ws.Range("J14").Value = 1.062418784727837
$ws.Range("K14").Value = 1.06001442542936
$ws.Range("L14").Value = 1.063384893104739
$ws.Range("M14").Value = 1.068473916002468
$ws.Range("N14").Value = 1.063927541263289

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.056231121310421
$ws.Range("D15").Value = 1.056706757171187
$ws.Range("E15").Value = 1.060135386203283
$ws.Range("F15").Value = 1.065224229746456
$ws.Range("I15").Value = 1.043172599218553
$ws.Range("J15").Value = 1.062621395806216
$ws.Range("K15").Value = 1.06018333237244
$ws.Range("L15").Value = 1.063599894878941
$ws.Range("M15").Value = 1.06867098671479
$ws.Range("N15").Value = 1.064130440072639

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.057626140503567
$ws.Range("D16").Value = 1.057801905919279
$ws.Range("E16").Value = 1.061499230424592
$ws.Range("F16").Value = 1.066482778054602
$ws.Range("I16").Value = 1.0435852020589
$ws.Range("J16").Value = 1.063797611370573
$ws.Range("K16").Value = 1.061163479264907
$ws.Range("L16").Value = 1.064848265529598
$ws.Range("M16").Value = 1.06981506418535
$ws.Range("N16").Value = 1.065308325998011

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.058498516056009
$ws.Range("D17").Value = 1.058486414044454
$ws.Range("E17").Value = 1.062352195875701
$ws.Range("F17").Value = 1.067269771385644
$ws.Range("I17").Value = 1.043842376526082
$ws.Range("J17").Value = 1.064532751243951
$ws.Range("K17").Value = 1.061775717312103
$ws.Range("L17").Value = 1.065628699357245
$ws.Range("M17").Value = 1.070530139440096
$ws.Range("N17").Value = 1.06604450985433

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.059006406149573
$ws.Range("D18").Value = 1.058884805055958
$ws.Range("E18").Value = 1.062848815986356
$ws.Range("F18").Value = 1.067727937694059
$ws.Range("I18").Value = 1.043991796153501
$ws.Range("J18").Value = 1.06496059712641
$ws.Range("K18").Value = 1.062131906106343
$ws.Range("L18").Value = 1.066082976858544
$ws.Range("M18").Value = 1.070946315552457
$ws.Range("N18").Value = 1.066472963327018

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.059179423520123
$ws.Range("D19").Value = 1.059020499467643
$ws.Range("E19").Value = 1.063017999270569
$ws.Range("F19").Value = 1.067884013627249
$ws.Range("I19").Value = 1.044042645502194
$ws.Range("J19").Value = 1.065106321815331
$ws.Range("K19").Value = 1.062253202406178
$ws.Range("L19").Value = 1.066237716203666
$ws.Range("M19").Value = 1.071088066794491
$ws.Range("N19").Value = 1.066618894961712

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.058405017204394
$ws.Range("D20").Value = 1.058413063157497
$ws.Range("E20").Value = 1.06226077419121
$ws.Range("F20").Value = 1.067185425104823
$ws.Range("I20").Value = 1.04381484485912
$ws.Range("J20").Value = 1.064453976122239
$ws.Range("K20").Value = 1.061710125304051
$ws.Range("L20").Value = 1.065545063350116
$ws.Range("M20").Value = 1.070453513575793
$ws.Range("N20").Value = 1.06596562286291

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.05587605423672
$ws.Range("D21").Value = 1.056427908733206
$ws.Range("E21").Value = 1.05978828049908
$ws.Range("F21").Value = 1.064903885890155
$ws.Range("I21").Value = 1.04306731975998
$ws.Range("J21").Value = 1.062321893977038
$ws.Range("K21").Value = 1.059933645160921
$ws.Range("L21").Value = 1.063282080949495
$ws.Range("M21").Value = 1.068379675167761
$ws.Range("N21").Value = 1.06383051291651

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.054277093432737
$ws.Range("D22").Value = 1.055171661004153
$ws.Range("E22").Value = 1.058225304208421
$ws.Range("F22").Value = 1.063461236172887
$ws.Range("I22").Value = 1.042591928862845
$ws.Range("J22").Value = 1.060972538951382
$ws.Range("K22").Value = 1.058808176340345
$ws.Range("L22").Value = 1.06185052405503
$ws.Range("M22").Value = 1.067067255058211
$ws.Range("N22").Value = 1.062479241651944

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.055125592574371
$ws.Range("D23").Value = 1.055838402235913
$ws.Range("E23").Value = 1.059054681032629
$ws.Range("F23").Value = 1.064226799118583
$ws.Range("I23").Value = 1.042844459279704
$ws.Range("J23").Value = 1.061688708471223
$ws.Range("K23").Value = 1.059405628682522
$ws.Range("L23").Value = 1.062610261399853
$ws.Range("M23").Value = 1.067763814244507
$ws.Range("N23").Value = 1.06319642821465

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.058447268249813
$ws.Range("D24").Value = 1.058446209955738
$ws.Range("E24").Value = 1.062302086494896
$ws.Range("F24").Value = 1.06722354025608
$ws.Range("I24").Value = 1.043827287043998
$ws.Range("J24").Value = 1.064489574138785
$ws.Range("K24").Value = 1.061739766348194
$ws.Range("L24").Value = 1.065582857753099
$ws.Range("M24").Value = 1.070488140331757
$ws.Range("N24").Value = 1.066001271432723

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.062264710741958
$ws.Range("D25").Value = 1.061438383489443
$ws.Range("E25").Value = 1.066035359290731
$ws.Range("F25").Value = 1.070666962439636
$ws.Range("I25").Value = 1.044944876486471
$ws.Range("J25").Value = 1.067702743317395
$ws.Range("K25").Value = 1.064412460128706
$ws.Range("L25").Value = 1.068995803815235
$ws.Range("M25").Value = 1.073613800958981
$ws.Range("N25").Value = 1.069219003680123
